$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add the new "formant values" worksheet after the last existing sheet
# (LCDViews), so it becomes the 4th / active tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "formant values"

# ---------------------------------------------------------------------
# Row 3 - scale constants used by the formula table below.
# ---------------------------------------------------------------------
$ws.Range("C3").Value = 9.5
$ws.Range("D3").Value = 19
$ws.Range("E3").Formula = "=256/36"

# ---------------------------------------------------------------------
# Row 6 - single label cell "i" above the little note-scale table.
# ---------------------------------------------------------------------
$ws.Range("D6").Value = "i"
$ws.Range("D6").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Rows 7-10 - frequency -> note-scale conversion table.
# ---------------------------------------------------------------------
$ws.Range("B7").Value = 450
$ws.Range("B8").Value = 2060
$ws.Range("B9").Value = 2700
$ws.Range("B10").Value = 3570

$ws.Range("C7").Formula = '=(B7+$C$3)/$D$3'
$ws.Range("C8:C10").Formula = '=(B8+$C$3)/$D$3'
$ws.Range("C7:C10").NumberFormat = "0"

$ws.Range("D7").Formula = "=C7/7"
$ws.Range("D8:D10").Formula = "=C8/7"
$ws.Range("D7:D10").NumberFormat = "0"
$ws.Range("D7:D10").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Rows 12-16 - "ee" / "i" comparison table.
# ---------------------------------------------------------------------
$ws.Range("B12").Value = "ee"
$ws.Range("C12").Value = "i"
$ws.Range("B12:C12").HorizontalAlignment = -4152

$ws.Range("B13").Value = 2.1766917293233083
$ws.Range("C13").Value = 3.4548872180451129
$ws.Range("B14").Value = 15.409774436090226
$ws.Range("C14").Value = 15.56015037593985
$ws.Range("B15").Value = 22.928571428571427
$ws.Range("C15").Value = 20.372180451127821
$ws.Range("B16").Value = 27.139097744360903
$ws.Range("C16").Value = 26.913533834586467

$ws.Range("B13:B16").NumberFormat = "0"
$ws.Range("B13:B16").HorizontalAlignment = -4152
$ws.Range("C13:C16").NumberFormat = "0"

# ---------------------------------------------------------------------
# Rows 25-29 - full formant value table for all vowel sounds.
# ---------------------------------------------------------------------
$ws.Range("G25").Value = "ee"
$ws.Range("H25").Value = "   I"
$ws.Range("I25").Value = "    e"
$ws.Range("J25").Value = "    a"
$ws.Range("K25").Value = "    o"
$ws.Range("L25").Value = "    oh"
$ws.Range("M25").Value = "   foot"
$ws.Range("N25").Value = " boot"
$ws.Range("O25").Value = " r"
$ws.Range("P25").Value = "    l"
$ws.Range("Q25").Value = "    uh"
$ws.Range("G25:Q25").HorizontalAlignment = -4152

$ws.Range("F26").Value = "f1"
$ws.Range("F26").HorizontalAlignment = -4152
$ws.Range("G26").Value = 280
$ws.Range("H26").Value = 450
$ws.Range("I26").Value = 550
$ws.Range("J26").Value = 700
$ws.Range("K26").Value = 775
$ws.Range("L26").Value = 575
$ws.Range("M26").Value = 425
$ws.Range("N26").Value = 275
$ws.Range("O26").Value = 560
$ws.Range("P26").Value = 560
$ws.Range("Q26").Value = 700

$ws.Range("F27").Value = "f2"
$ws.Range("F27").HorizontalAlignment = -4152
$ws.Range("G27").Value = 2040
$ws.Range("H27").Value = 2060
$ws.Range("I27").Value = 1950
$ws.Range("J27").Value = 1800
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = 1000
$ws.Range("N27").Value = 850
$ws.Range("O27").Value = 1200
$ws.Range("P27").Value = 820
$ws.Range("Q27").Value = 1300

$ws.Range("F28").Value = "f3"
$ws.Range("F28").HorizontalAlignment = -4152
$ws.Range("G28").Value = 3040
$ws.Range("H28").Value = 2700
$ws.Range("I28").Value = 2600
$ws.Range("J28").Value = 2550
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 2450
$ws.Range("M28").Value = 2400
$ws.Range("N28").Value = 2400
$ws.Range("O28").Value = 1500
$ws.Range("P28").Value = 2700
$ws.Range("Q28").Value = 2600

$ws.Range("F29").Value = "f4"
$ws.Range("F29").HorizontalAlignment = -4152
$ws.Range("G29").Value = 3600
$ws.Range("H29").Value = 3570
$ws.Range("I29").Value = 3400
$ws.Range("J29").Value = 3400
$ws.Range("K29").Value = 3500
$ws.Range("L29").Value = 3500
$ws.Range("M29").Value = 3500
$ws.Range("N29").Value = 3500
$ws.Range("O29").Value = 3050
$ws.Range("P29").Value = 3600
$ws.Range("Q29").Value = 3100

# ---------------------------------------------------------------------
# Column C width (values are best-fit in the original workbook) and the
# selection left active on the new sheet.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 8.6

$ws.Range("B12:C16").Select()
